$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 119,3
$data[0,0] = "SRR8423796"
$data[0,1] = 1624
$data[0,2] = "Illumina MiSeq"
$data[1,0] = "SRR8423797"
$data[1,1] = 108455
$data[1,2] = "Illumina MiSeq"
$data[2,0] = "SRR8423798"
$data[2,1] = 97249
$data[2,2] = "Illumina MiSeq"
$data[3,0] = "SRR8423799"
$data[3,1] = 1653
$data[3,2] = "Illumina MiSeq"
$data[4,0] = "SRR8423800"
$data[4,1] = 258794
$data[4,2] = "Illumina MiSeq"
$data[5,0] = "SRR8423801"
$data[5,1] = 121846
$data[5,2] = "Illumina MiSeq"
$data[6,0] = "SRR8423802"
$data[6,1] = 211903
$data[6,2] = "Illumina MiSeq"
$data[7,0] = "SRR8423803"
$data[7,1] = 261720
$data[7,2] = "Illumina MiSeq"
$data[8,0] = "SRR8423804"
$data[8,1] = 114730
$data[8,2] = "Illumina MiSeq"
$data[9,0] = "SRR8423805"
$data[9,1] = 119611
$data[9,2] = "Illumina MiSeq"
$data[10,0] = "SRR8423806"
$data[10,1] = 1459396
$data[10,2] = "Illumina MiSeq"
$data[11,0] = "SRR8423807"
$data[11,1] = 951233
$data[11,2] = "Illumina MiSeq"
$data[12,0] = "SRR8423808"
$data[12,1] = 1490191
$data[12,2] = "Illumina MiSeq"
$data[13,0] = "SRR8423809"
$data[13,1] = 1369757
$data[13,2] = "Illumina MiSeq"
$data[14,0] = "SRR8423810"
$data[14,1] = 147238
$data[14,2] = "Illumina MiSeq"
$data[15,0] = "SRR8423811"
$data[15,1] = 143700
$data[15,2] = "Illumina MiSeq"
$data[16,0] = "SRR8423812"
$data[16,1] = 29323
$data[16,2] = "Illumina MiSeq"
$data[17,0] = "SRR8423813"
$data[17,1] = 185108
$data[17,2] = "Illumina MiSeq"
$data[18,0] = "SRR8423814"
$data[18,1] = 1613157
$data[18,2] = "Illumina MiSeq"
$data[19,0] = "SRR8423816"
$data[19,1] = 13537348
$data[19,2] = "Illumina NovaSeq 6000"
$data[20,0] = "SRR8423817"
$data[20,1] = 11305252
$data[20,2] = "Illumina NovaSeq 6000"
$data[21,0] = "SRR8423818"
$data[21,1] = 9873015
$data[21,2] = "Illumina NovaSeq 6000"
$data[22,0] = "SRR8423819"
$data[22,1] = 9374518
$data[22,2] = "Illumina NovaSeq 6000"
$data[23,0] = "SRR8423820"
$data[23,1] = 15399680
$data[23,2] = "Illumina NovaSeq 6000"
$data[24,0] = "SRR8423821"
$data[24,1] = 19915070
$data[24,2] = "Illumina NovaSeq 6000"
$data[25,0] = "SRR8423822"
$data[25,1] = 11039286
$data[25,2] = "Illumina NovaSeq 6000"
$data[26,0] = "SRR8423823"
$data[26,1] = 12672547
$data[26,2] = "Illumina NovaSeq 6000"
$data[27,0] = "SRR8423824"
$data[27,1] = 11986897
$data[27,2] = "Illumina NovaSeq 6000"
$data[28,0] = "SRR8423825"
$data[28,1] = 14836573
$data[28,2] = "Illumina NovaSeq 6000"
$data[29,0] = "SRR8423826"
$data[29,1] = 9833029
$data[29,2] = "Illumina NovaSeq 6000"
$data[30,0] = "SRR8423827"
$data[30,1] = 13488303
$data[30,2] = "Illumina NovaSeq 6000"
$data[31,0] = "SRR8423828"
$data[31,1] = 20515178
$data[31,2] = "Illumina NovaSeq 6000"
$data[32,0] = "SRR8423829"
$data[32,1] = 12975096
$data[32,2] = "Illumina NovaSeq 6000"
$data[33,0] = "SRR8423830"
$data[33,1] = 13794276
$data[33,2] = "Illumina NovaSeq 6000"
$data[34,0] = "SRR8423831"
$data[34,1] = 18010530
$data[34,2] = "Illumina NovaSeq 6000"
$data[35,0] = "SRR8423832"
$data[35,1] = 28063231
$data[35,2] = "Illumina NovaSeq 6000"
$data[36,0] = "SRR8423833"
$data[36,1] = 26039122
$data[36,2] = "Illumina NovaSeq 6000"
$data[37,0] = "SRR8423834"
$data[37,1] = 17720066
$data[37,2] = "Illumina NovaSeq 6000"
$data[38,0] = "SRR8423835"
$data[38,1] = 18898195
$data[38,2] = "Illumina NovaSeq 6000"
$data[39,0] = "SRR8423836"
$data[39,1] = 14027792
$data[39,2] = "Illumina NovaSeq 6000"
$data[40,0] = "SRR8423837"
$data[40,1] = 20016455
$data[40,2] = "Illumina NovaSeq 6000"
$data[41,0] = "SRR8423838"
$data[41,1] = 104614
$data[41,2] = "Illumina MiSeq"
$data[42,0] = "SRR8423839"
$data[42,1] = 119558
$data[42,2] = "Illumina MiSeq"
$data[43,0] = "SRR8423840"
$data[43,1] = 12264113
$data[43,2] = "Illumina NovaSeq 6000"
$data[44,0] = "SRR8423841"
$data[44,1] = 9539568
$data[44,2] = "Illumina NovaSeq 6000"
$data[45,0] = "SRR8423842"
$data[45,1] = 11595236
$data[45,2] = "Illumina NovaSeq 6000"
$data[46,0] = "SRR8423843"
$data[46,1] = 14272314
$data[46,2] = "Illumina NovaSeq 6000"
$data[47,0] = "SRR8423844"
$data[47,1] = 122979
$data[47,2] = "Illumina MiSeq"
$data[48,0] = "SRR8423845"
$data[48,1] = 12587
$data[48,2] = "Illumina MiSeq"
$data[49,0] = "SRR8423846"
$data[49,1] = 13031727
$data[49,2] = "Illumina NovaSeq 6000"
$data[50,0] = "SRR8423847"
$data[50,1] = 300981
$data[50,2] = "Illumina NovaSeq 6000"
$data[51,0] = "SRR8423848"
$data[51,1] = 921553
$data[51,2] = "Illumina NovaSeq 6000"
$data[52,0] = "SRR8423849"
$data[52,1] = 10678846
$data[52,2] = "Illumina NovaSeq 6000"
$data[53,0] = "SRR8423850"
$data[53,1] = 15955302
$data[53,2] = "Illumina NovaSeq 6000"
$data[54,0] = "SRR8423851"
$data[54,1] = 16187888
$data[54,2] = "Illumina NovaSeq 6000"
$data[55,0] = "SRR8423852"
$data[55,1] = 15571061
$data[55,2] = "Illumina NovaSeq 6000"
$data[56,0] = "SRR8423853"
$data[56,1] = 14855966
$data[56,2] = "Illumina NovaSeq 6000"
$data[57,0] = "SRR8423854"
$data[57,1] = 466603
$data[57,2] = "Illumina NovaSeq 6000"
$data[58,0] = "SRR8423855"
$data[58,1] = 14245489
$data[58,2] = "Illumina NovaSeq 6000"
$data[59,0] = "SRR8423856"
$data[59,1] = 115191
$data[59,2] = "Illumina MiSeq"
$data[60,0] = "SRR8423857"
$data[60,1] = 201602
$data[60,2] = "Illumina MiSeq"
$data[61,0] = "SRR8423859"
$data[61,1] = 97890
$data[61,2] = "Illumina MiSeq"
$data[62,0] = "SRR8423860"
$data[62,1] = 133854
$data[62,2] = "Illumina MiSeq"
$data[63,0] = "SRR8423862"
$data[63,1] = 161675
$data[63,2] = "Illumina MiSeq"
$data[64,0] = "SRR8423863"
$data[64,1] = 149809
$data[64,2] = "Illumina MiSeq"
$data[65,0] = "SRR8423864"
$data[65,1] = 208897
$data[65,2] = "Illumina MiSeq"
$data[66,0] = "SRR8423865"
$data[66,1] = 216386
$data[66,2] = "Illumina MiSeq"
$data[67,0] = "SRR8423866"
$data[67,1] = 13456413
$data[67,2] = "Illumina NovaSeq 6000"
$data[68,0] = "SRR8423867"
$data[68,1] = 14179351
$data[68,2] = "Illumina NovaSeq 6000"
$data[69,0] = "SRR8423868"
$data[69,1] = 194866
$data[69,2] = "Illumina MiSeq"
$data[70,0] = "SRR8423869"
$data[70,1] = 718
$data[70,2] = "Illumina MiSeq"
$data[71,0] = "SRR8423870"
$data[71,1] = 20863983
$data[71,2] = "Illumina NovaSeq 6000"
$data[72,0] = "SRR8423871"
$data[72,1] = 17109238
$data[72,2] = "Illumina NovaSeq 6000"
$data[73,0] = "SRR8423872"
$data[73,1] = 14084641
$data[73,2] = "Illumina NovaSeq 6000"
$data[74,0] = "SRR8423873"
$data[74,1] = 16305228
$data[74,2] = "Illumina NovaSeq 6000"
$data[75,0] = "SRR8423874"
$data[75,1] = 17899271
$data[75,2] = "Illumina NovaSeq 6000"
$data[76,0] = "SRR8423875"
$data[76,1] = 22275436
$data[76,2] = "Illumina NovaSeq 6000"
$data[77,0] = "SRR8423876"
$data[77,1] = 1326845
$data[77,2] = "Illumina MiSeq"
$data[78,0] = "SRR8423877"
$data[78,1] = 764657
$data[78,2] = "Illumina MiSeq"
$data[79,0] = "SRR8423878"
$data[79,1] = 21
$data[79,2] = "Illumina MiSeq"
$data[80,0] = "SRR8423879"
$data[80,1] = 175017
$data[80,2] = "Illumina MiSeq"
$data[81,0] = "SRR8423880"
$data[81,1] = 145708
$data[81,2] = "Illumina MiSeq"
$data[82,0] = "SRR8423881"
$data[82,1] = 139661
$data[82,2] = "Illumina MiSeq"
$data[83,0] = "SRR8423882"
$data[83,1] = 20240060
$data[83,2] = "Illumina NovaSeq 6000"
$data[84,0] = "SRR8423883"
$data[84,1] = 6845865
$data[84,2] = "Illumina NovaSeq 6000"
$data[85,0] = "SRR8423884"
$data[85,1] = 13668744
$data[85,2] = "Illumina NovaSeq 6000"
$data[86,0] = "SRR8423885"
$data[86,1] = 14147968
$data[86,2] = "Illumina NovaSeq 6000"
$data[87,0] = "SRR8423886"
$data[87,1] = 183712
$data[87,2] = "Illumina MiSeq"
$data[88,0] = "SRR8423887"
$data[88,1] = 175829
$data[88,2] = "Illumina MiSeq"
$data[89,0] = "SRR8423888"
$data[89,1] = 172378
$data[89,2] = "Illumina MiSeq"
$data[90,0] = "SRR8423889"
$data[90,1] = 187223
$data[90,2] = "Illumina MiSeq"
$data[91,0] = "SRR8423890"
$data[91,1] = 108926
$data[91,2] = "Illumina MiSeq"
$data[92,0] = "SRR8423891"
$data[92,1] = 146692
$data[92,2] = "Illumina MiSeq"
$data[93,0] = "SRR8423892"
$data[93,1] = 158588
$data[93,2] = "Illumina MiSeq"
$data[94,0] = "SRR8423893"
$data[94,1] = 145661
$data[94,2] = "Illumina MiSeq"
$data[95,0] = "SRR8423894"
$data[95,1] = 186246
$data[95,2] = "Illumina MiSeq"
$data[96,0] = "SRR8423895"
$data[96,1] = 119946
$data[96,2] = "Illumina MiSeq"
$data[97,0] = "SRR8423896"
$data[97,1] = 205318
$data[97,2] = "Illumina MiSeq"
$data[98,0] = "SRR8423897"
$data[98,1] = 184862
$data[98,2] = "Illumina MiSeq"
$data[99,0] = "SRR8423898"
$data[99,1] = 1112965
$data[99,2] = "Illumina MiSeq"
$data[100,0] = "SRR8423899"
$data[100,1] = 972280
$data[100,2] = "Illumina MiSeq"
$data[101,0] = "SRR8423900"
$data[101,1] = 984088
$data[101,2] = "Illumina MiSeq"
$data[102,0] = "SRR8423901"
$data[102,1] = 888371
$data[102,2] = "Illumina MiSeq"
$data[103,0] = "SRR8423902"
$data[103,1] = 1279222
$data[103,2] = "Illumina MiSeq"
$data[104,0] = "SRR8423903"
$data[104,1] = 1050684
$data[104,2] = "Illumina MiSeq"
$data[105,0] = "SRR8423904"
$data[105,1] = 1272635
$data[105,2] = "Illumina MiSeq"
$data[106,0] = "SRR8423905"
$data[106,1] = 1132226
$data[106,2] = "Illumina MiSeq"
$data[107,0] = "SRR8423906"
$data[107,1] = 1174312
$data[107,2] = "Illumina MiSeq"
$data[108,0] = "SRR8423907"
$data[108,1] = 1197633
$data[108,2] = "Illumina MiSeq"
$data[109,0] = "SRR8423908"
$data[109,1] = 96533
$data[109,2] = "Illumina MiSeq"
$data[110,0] = "SRR8423909"
$data[110,1] = 169524
$data[110,2] = "Illumina MiSeq"
$data[111,0] = "SRR8423910"
$data[111,1] = 136907
$data[111,2] = "Illumina MiSeq"
$data[112,0] = "SRR8423911"
$data[112,1] = 173420
$data[112,2] = "Illumina MiSeq"
$data[113,0] = "SRR8423912"
$data[113,1] = 173365
$data[113,2] = "Illumina MiSeq"
$data[114,0] = "SRR8423913"
$data[114,1] = 130294
$data[114,2] = "Illumina MiSeq"
$data[115,0] = "SRR8423914"
$data[115,1] = 156235
$data[115,2] = "Illumina MiSeq"
$data[116,0] = "SRR8423915"
$data[116,1] = 131060
$data[116,2] = "Illumina MiSeq"
$data[117,0] = "SRR8423916"
$data[117,1] = 169550
$data[117,2] = "Illumina MiSeq"
$data[118,0] = "SRR8423917"
$data[118,1] = 185600
$data[118,2] = "Illumina MiSeq"

$ws.Range("A2:C120").Value = $data

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A120"))
$ws.Sort.SetRange($ws.Range("A2:D120"))
$ws.Sort.Header = -4163
$ws.Sort.Apply()

$ws.Range("A2").Select()
